$d = $word.ActiveDocument

# Prevent Word's smart-quote autocorrect from mangling straight quotes
# that merely pass through the Find/Replace on their way to being saved.
try { $word.Options.AutoFormatAsYouTypeReplaceQuotes = $false } catch {}
try { $word.Options.AutoFormatReplaceQuotes = $false } catch {}

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false,
                             $true, 1, $false, $replace, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 1. Intro paragraph: date/time, lessor/lessee names, locations
#    (split around the straight-quote marks ' ... ' so Find/Replace
#    never re-types them - otherwise Word's smart-quote autocorrect
#    turns them into curly quotes even with AutoFormatReplaceQuotes
#    disabled)
# ---------------------------------------------------------------------
Replace-Text "This Deed of Lease is made at 20:30 this Monday day of September, 2024 between Mihir of Maharashtra hereinafter called " "This Deed of Lease is made at 12:03 this 1 day of march, 2025 between sathvik of hyderabad hereinafter called "
Replace-Text " of the One Part and Raju also of Gujarat hereinafter called " " of the One Part and bandi also of hyderabad hereinafter called "

# ---------------------------------------------------------------------
# 2. Clause 1: location, commencing date, rent amount, payment date
# ---------------------------------------------------------------------
Replace-Text "situated at dfdf and described" "situated at hyderabad and described"
Replace-Text "commencing from the 1st day of September, 2024, but subject" "commencing from the 1st day of march, 1, but subject"
Replace-Text "monthly ground rent of Rs 12000 free and clear" "monthly ground rent of Rs 10000 free and clear"
Replace-Text "shall be paid on the 5th day of September and the subsequent" "shall be paid on the 5th day of march and the subsequent"

# ---------------------------------------------------------------------
# 3. Clause a: payment date, interest rate
# ---------------------------------------------------------------------
Replace-Text "shall be paid on the 5th of September and the subsequent" "shall be paid on the 5th of march and the subsequent"
Replace-Text "at the rate of 2 % per annum" "at the rate of 10 % per annum"

# ---------------------------------------------------------------------
# 4. Clause 4: arrears months, notice months
# ---------------------------------------------------------------------
Replace-Text "for the space of 3 months after" "for the space of 2 months after"
Replace-Text "carry out the same within 3 months from" "carry out the same within 2 months from"

# ---------------------------------------------------------------------
# 5. Signed and delivered - Lessor paragraph (single run)
# ---------------------------------------------------------------------
Replace-Text "Withinnamed Lessor Raju in the presence of dfdnf" "Withinnamed Lessor bandi in the presence of suii"

# ---------------------------------------------------------------------
# 6. Signed and delivered - Lessee paragraph.
#    This paragraph is made of THREE runs that all share identical
#    character formatting: "Withinnamed Lessee ", "Mihir",
#    " in the presence of dfdjf". Only the last two runs change text
#    (Mihir -> sathvik, dfdjf -> poppy); the first run is untouched.
#    A plain Find/Replace (or any Range.Text edit) on this paragraph
#    causes same-formatted adjacent runs to coalesce into one run on
#    save, so we briefly toggle a binary formatting flag (Italic) on
#    the neighbouring runs before editing - this keeps them from
#    merging with the run being edited - then toggle it back off.
# ---------------------------------------------------------------------
$lesseeText = "Withinnamed Lessee Mihir in the presence of dfdjf"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $lesseeText) {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $full = $target.Range.Text

    $idxName = $full.IndexOf("Mihir")
    $nameStart = $pStart + $idxName
    $nameEnd = $nameStart + 5

    $idxWit = $full.IndexOf(" in the presence of dfdjf")
    $witStart = $pStart + $idxWit
    $witEnd = $witStart + 26

    # Separate run 1 ("Withinnamed Lessee ") and run 3 (" in the
    # presence of dfdjf") from run 2 ("Mihir") while we edit, so they
    # don't get coalesced together.
    $run1 = $d.Range($pStart, $nameStart)
    $run1.Italic = 1

    $run3 = $d.Range($witStart, $witEnd)
    $run3.Italic = 1

    # Edit run 2: Mihir -> sathvik
    $run2 = $d.Range($nameStart, $nameEnd)
    $run2.Text = "sathvik"

    # Recompute run 3's offsets (run 2 grew by 2 characters)
    $full2 = $target.Range.Text
    $idxWit2 = $full2.IndexOf(" in the presence of dfdjf")
    $witStart2 = $pStart + $idxWit2
    $witEnd2 = $witStart2 + 26

    # Edit run 3: " in the presence of dfdjf" -> " in the presence of poppy"
    $run3b = $d.Range($witStart2, $witEnd2)
    $run3b.Text = " in the presence of poppy"

    # Restore formatting (turn Italic back off) on run 1 and run 3
    $full3 = $target.Range.Text
    $idxName3 = $full3.IndexOf("sathvik")
    $nameStart3 = $pStart + $idxName3
    $run1b = $d.Range($pStart, $nameStart3)
    $run1b.Italic = 0

    $idxWit3 = $full3.IndexOf(" in the presence of poppy")
    $witStart3 = $pStart + $idxWit3
    $witEnd3 = $witStart3 + 26
    $run3c = $d.Range($witStart3, $witEnd3)
    $run3c.Italic = 0
}
